{"js": "const replacements = [\n  [\"2024-01-10 Wednesday\", \"2024-01-11 Thursday\"],\n\n  [\"410\u00f72=205, 0\", \"926\u00f77=132, 2\"],\n  [\"168\u00f73=56, 0\", \"486\u00f74=121, 2\"],\n  [\"682\u00f76=113, 4\", \"873\u00f79=97, 0\"],\n  [\"956\u00f72=478, 0\", \"696\u00f73=232, 0\"],\n  [\"259\u00f73=86, 1\", \"269\u00f78=33, 5\"],\n\n  [\"414\u00f75=82, 4\", \"179\u00f72=89, 1\"],\n  [\"955\u00f75=191, 0\", \"831\u00f77=118, 5\"],\n  [\"333\u00f72=166, 1\", \"723\u00f75=144, 3\"],\n  [\"284\u00f79=31, 5\", \"360\u00f72=180, 0\"],\n  [\"585\u00f75=117, 0\", \"987\u00f78=123, 3\"],\n\n  [\"744\u00f73=248, 0\", \"539\u00f77=77, 0\"],\n  [\"510\u00f76=85, 0\", \"778\u00f78=97, 2\"],\n  [\"404\u00f74=101, 0\", \"118\u00f78=14, 6\"],\n  [\"555\u00f78=69, 3\", \"619\u00f78=77, 3\"],\n  [\"927\u00f77=132, 3\", \"746\u00f76=124, 2\"],\n\n  [\"893\u00f73=297, 2\", \"397\u00f76=66, 1\"],\n  [\"677\u00f74=169, 1\", \"350\u00f78=43, 6\"],\n  [\"207\u00f75=41, 2\", \"222\u00f79=24, 6\"],\n  [\"485\u00f76=80, 5\", \"632\u00f73=210, 2\"],\n  [\"340\u00f79=37, 7\", \"884\u00f75=176, 4\"],\n\n  [\"479\u00f73=159, 2\", \"852\u00f72=426, 0\"],\n  [\"932\u00f74=233, 0\", \"650\u00f75=130, 0\"],\n  [\"101\u00f79=11, 2\", \"224\u00f77=32, 0\"],\n  [\"377\u00f74=94, 1\", \"980\u00f77=140, 0\"],\n  [\"620\u00f73=206, 2\", \"265\u00f72=132, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\nReplace-Text \"2024-01-10 Wednesday\" \"2024-01-11 Thursday\"\n\nReplace-Text \"410\u00f72=205, 0\" \"926\u00f77=132, 2\"\nReplace-Text \"168\u00f73=56, 0\" \"486\u00f74=121, 2\"\nReplace-Text \"682\u00f76=113, 4\" \"873\u00f79=97, 0\"\nReplace-Text \"956\u00f72=478, 0\" \"696\u00f73=232, 0\"\nReplace-Text \"259\u00f73=86, 1\" \"269\u00f78=33, 5\"\n\nReplace-Text \"414\u00f75=82, 4\" \"179\u00f72=89, 1\"\nReplace-Text \"955\u00f75=191, 0\" \"831\u00f77=118, 5\"\nReplace-Text \"333\u00f72=166, 1\" \"723\u00f75=144, 3\"\nReplace-Text \"284\u00f79=31, 5\" \"360\u00f72=180, 0\"\nReplace-Text \"585\u00f75=117, 0\" \"987\u00f78=123, 3\"\n\nReplace-Text \"744\u00f73=248, 0\" \"539\u00f77=77, 0\"\nReplace-Text \"510\u00f76=85, 0\" \"778\u00f78=97, 2\"\nReplace-Text \"404\u00f74=101, 0\" \"118\u00f78=14, 6\"\nReplace-Text \"555\u00f78=69, 3\" \"619\u00f78=77, 3\"\nReplace-Text \"927\u00f77=132, 3\" \"746\u00f76=124, 2\"\n\nReplace-Text \"893\u00f73=297, 2\" \"397\u00f76=66, 1\"\nReplace-Text \"677\u00f74=169, 1\" \"350\u00f78=43, 6\"\nReplace-Text \"207\u00f75=41, 2\" \"222\u00f79=24, 6\"\nReplace-Text \"485\u00f76=80, 5\" \"632\u00f73=210, 2\"\nReplace-Text \"340\u00f79=37, 7\" \"884\u00f75=176, 4\"\n\nReplace-Text \"479\u00f73=159, 2\" \"852\u00f72=426, 0\"\nReplace-Text \"932\u00f74=233, 0\" \"650\u00f75=130, 0\"\nReplace-Text \"101\u00f79=11, 2\" \"224\u00f77=32, 0\"\nReplace-Text \"377\u00f74=94, 1\" \"980\u00f77=140, 0\"\nReplace-Text \"620\u00f73=206, 2\" \"265\u00f72=132, 1\"\n"}
